# Hangulatvilágítás.xlsx — add the "Soros protokoll" sheet documenting the
# basic serial communication protocol with the PC, and move the selection
# on the first sheet down to where the interpolation section now lives.

$wb = $excel.ActiveWorkbook

# --- 1. Update the selection / scroll position on the first sheet --------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D24").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 16
    $win.ScrollColumn = 1
}

# --- 2. Add the new "Soros protokoll" worksheet as the last tab ----------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Soros protokoll"

# --- 3. Fill in the protocol table ---------------------------------------
$ws3.Range("A1").Value = "Üzenettípus"
$ws3.Range("B1").Value = "Formátum"
$ws3.Range("C1").Value = "Tartalom"

$ws3.Range("A2").Value = "Fejléc"
$ws3.Range("B2").Value = """HEADER [headerbytes]"""
$ws3.Range("C2").Value = "Bájtfolyamként a programheader"

$ws3.Range("A3").Value = "Kvantum"
$ws3.Range("B3").Value = """Q [quantumbytes]"
$ws3.Range("C3").Value = "Bájtfolyamként a kvantum"

$ws3.Range("A4").Value = "Start"
$ws3.Range("B4").Value = """START"""
$ws3.Range("C4").Value = "Új program leküldésének inicializálása"

$ws3.Range("A5").Value = "Stop"
$ws3.Range("B5").Value = """STOP"""
$ws3.Range("C5").Value = "Programleküldés vége"

$ws3.Range("A6").Value = "Futtatás"
$ws3.Range("B6").Value = """RUN [ID]"""
$ws3.Range("C6").Value = "Szövegként a futtatandó program ID-je"

# --- 4. Column widths (bestFit-ish widths matching the authored sheet) ---
$ws3.Columns.Item(1).ColumnWidth = 12.833333333333332
$ws3.Columns.Item(2).ColumnWidth = 20.5
$ws3.Columns.Item(3).ColumnWidth = 33.833333333333336

# --- 5. Leave the selection where the author left it, just below the table
$ws3.Range("B7").Select()
